$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Week 2 (column E) estimates for sprint backlog tasks (rows 4-16 and 26-36)
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 2
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 2
$ws.Range("E10").Value = 3
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 3

$ws.Range("E26").Value = 2
$ws.Range("E27").Value = 2
$ws.Range("E28").Value = 2
$ws.Range("E29").Value = 2
$ws.Range("E30").Value = 3
$ws.Range("E31").Value = 3
$ws.Range("E32").Value = 2
$ws.Range("E33").Value = 2
$ws.Range("E34").Value = 2
$ws.Range("E35").Value = 3
$ws.Range("E36").Value = 3

# Refresh the burndown chart so its cached series values (Week 2 / column E,
# row 44 total) reflect the updated totals
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.Refresh()

# Update the active sheet view position / selection to match the saved state
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("G17").Select()
